$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A9").Value = -1
$ws.Range("A9").Select()
